# Updated cryptos list on Tue Nov 19 07:54:05 UTC 2024 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for the
# crypto rows, and swaps row 51 from VeChain to OKB.
#
# Several Price values look like plain decimals (e.g. "243.69"); a bare
# assignment would let Excel auto-convert those to numbers, which would
# change the cell's stored type away from the original text. Prefixing
# those with a leading apostrophe keeps them as literal text (the
# apostrophe itself is just a text-entry marker and is not stored as part
# of the value), matching the original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.836.65'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '3.123.59'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'243.69"
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = "'618.07"
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('D8').Value = "'0.391"
$ws.Range('E8').Value = '  +4.11%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.119.32'
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').Value = "'0.757"
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = "'0.0000253"
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').Value = "'35.18"
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('D15').Value = "'5.60"
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').Value = '91.679.08'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D18').Value = '3.146.03'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = "'14.93"
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').Value = "'5.87"
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = "'455.62"
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('E23').Value = '  -5.01%  '
$ws.Range('D24').Value = "'9.18"
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = "'5.90"
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').Value = "'89.62"
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('E27').Value = '  +46.67%  '
$ws.Range('E28').Value = '  -3.23%  '
$ws.Range('D30').Value = "'0.998"
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').Value = "'0.139"
$ws.Range('E31').Value = '  +16.88%  '
$ws.Range('D32').Value = "'0.227"
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E33').Value = '  -7.07%  '
$ws.Range('D34').Value = "'9.33"
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').Value = "'0.173"
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D36').Value = "'2.09"
$ws.Range('E36').Value = '  +7.91%  '
$ws.Range('D37').Value = "'26.34"
$ws.Range('E37').Value = '  -2.37%  '
$ws.Range('D38').Value = "'7.52"
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').Value = "'493.02"
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('D40').Value = "'1.31"
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('D41').Value = "'3.85"
$ws.Range('E41').Value = '  -7.99%  '
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('D43').Value = "'3.43"
$ws.Range('E43').Value = '  -6.00%  '
$ws.Range('D44').Value = "'22.20"
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D46').Value = "'159.04"
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = "'44.01"
$ws.Range('E51').Value = '  -2.10%  '
